$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sloth's "Activité 11" entry was wrongly marked "Non rendu" -- it was actually
# received. Correct the value and give it its own "Reçu" look (plain font,
# new lighter fill) instead of the "Non rendu" grey-italic style.
$rc = $ws.Range("E2")
$rc.Value = "Reçu"
$rc.Font.Italic = $false
$rc.Font.ThemeColor = 1
$rc.Interior.Color = 8830443

# Skido's row keeps displaying "Non rendu" with its original look; nothing
# visually changes there, only internal bookkeeping shifts to account for the
# newly introduced style/string.
